$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Copy the formatting of the last filled row (45) down onto row 46 so the
#     new row matches the established look (borders, number formats, fonts). ---
$ws.Range("C45:L45").Copy()
$ws.Range("C46:L46").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Give row 46 the same visual row height as the other filled rows.
$ws.Rows.Item(46).RowHeight = 15

# --- Fill in the newly-practiced test row (row 46) ---
$ws.Range("C46").Value = 42
$ws.Range("D46").Value = 45530
$ws.Range("E46").Value = "Official IELTS Practice material 2"
$ws.Range("F46").Value = 31
$ws.Range("G46").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"
$ws.Range("H46").Value = 31
$ws.Range("I46").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"
$ws.Range("J46").Value = 6
$ws.Range("K46").Value = 4
$ws.Range("L46").Formula = "=(G46+I46+J46+K46)/4"

# --- Column E got a bit wider to fit the longer course name. ---
$ws.Columns.Item(5).ColumnWidth = 26.5

# --- Selection moved on to the next empty row while the user kept working. ---
$ws.Range("K47").Select()
